$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update countries (Pais column) ---
# Cuba's case counts overtook Bulgaria's, so the two countries swap rows
# (the sorted-by-total-cases list re-ranks them).
$ws.Range("A79").Value = "Cuba"
$ws.Range("A80").Value = "Bulgaria"

# Namibia and San Vicente y las Granadinas are tied, and swap rows too.
$ws.Range("A194").Value = "Namibia"
$ws.Range("A195").Value = "San Vicente y las Granadinas"

# --- Update provincias / country stats ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1167264
$ws.Range("C4").Value = 6490
$ws.Range("D4").Value = 174017
$ws.Range("E4").Value = 925642
$ws.Range("G4").Value = 161
$ws.Range("H4").Value = 67605

# Row 27
$ws.Range("B27").Value = 20084
$ws.Range("C27").Value = 1062
$ws.Range("D27").Value = 5114
$ws.Range("E27").Value = 14513
$ws.Range("G27").Value = 20
$ws.Range("H27").Value = 457

# Row 29
$ws.Range("D29").Value = 1408
$ws.Range("E29").Value = 16779
$ws.Range("F29").Value = 22
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 18

# Row 61
$ws.Range("B61").Value = 3824
$ws.Range("C61").Value = 12
$ws.Range("D61").Value = 3379
$ws.Range("E61").Value = 349
$ws.Range("F61").Value = 22
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 96

# Row 79 (now Cuba)
$ws.Range("B79").Value = 1649
$ws.Range("C79").Value = 38
$ws.Range("D79").Value = 827
$ws.Range("E79").Value = 755
$ws.Range("F79").Value = 13
$ws.Range("H79").Value = 67

# Row 80 (now Bulgaria)
$ws.Range("B80").Value = 1618
$ws.Range("C80").Value = 24
$ws.Range("D80").Value = 308
$ws.Range("E80").Value = 1237
$ws.Range("F80").Value = 39
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 73

# Row 124
$ws.Range("D124").Value = 77
$ws.Range("E124").Value = 274
